$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.995.74"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "1.756.15"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "336.15"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9988"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3421"
$ws.Range("E8").Value = "  +0.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.07"
$ws.Range("E9").Value = "  -2.42%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.124"
$ws.Range("E10").Value = "  -1.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07237"
$ws.Range("E11").Value = "  -2.15%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.62"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.174"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.157"
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "1.751.83"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001061"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06608"
$ws.Range("E18").Value = "  -0.58%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "79.34"
$ws.Range("E19").Value = "  -3.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9992"
$ws.Range("E20").Value = "  -0.43%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.75"
$ws.Range("E21").Value = "  -3.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.222"
$ws.Range("E22").Value = "  -2.72%  "
$ws.Range("D23").Value = "27.983.23"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("E24").Value = "  -3.06%  "
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.89"
$ws.Range("E26").Value = "  +0.49%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.86"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.315"
$ws.Range("E28").Value = "  -5.01%  "
$ws.Range("D29").Value = "1.952.56"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("E30").Value = "  -10.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "131.55"
$ws.Range("E31").Value = "  -2.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.025"
$ws.Range("E32").Value = "  +1.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.853"
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08820"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.23"
$ws.Range("E35").Value = "  -3.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6580"
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02296"
$ws.Range("E37").Value = "  -4.55%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.164"
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06184"
$ws.Range("E39").Value = "  -1.26%  "
$ws.Range("B40").Value = "WEMIXTOKEN"
$ws.Range("C40").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.510"
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2109"
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.215"
$ws.Range("E42").Value = "  -3.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.968"
$ws.Range("E43").Value = "  -4.04%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9988"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.78"
$ws.Range("E45").Value = "  -3.38%  "
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.838"
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6073"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.09"
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.009"
$ws.Range("E49").Value = "  -3.24%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.173"
$ws.Range("E50").Value = "  +2.60%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.116"
$ws.Range("E51").Value = "  +5.25%  "
